$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update wording of the "isbtnAnnulerValider" explanation (column B, row 13)
$ws.Range("B13").Value = "Si le statut est intermédiaire ( >SAISI et <Final) et que l'utilisateur a effectué l'action de validation précédente et qu'aucun des supérieurs n'a effectué l'une des actions n++"

# Update wording of the "isbtnClose" explanation (column B, row 6)
$ws.Range("B6").Value = "Tant qu'on n'a pas validé final, si l'entité n'est pas close et l'utilisateur est celui qui a saisi alors il peut close"

# Update the view: scroll back to top and move the active selection to B24
$ws.Range("A1").Select()
$ws.Range("B24").Select()
